# Change wording of "QR Code" to "pass" on the online check-in slide
# ("Print Your QR Code & bring it with you" -> "Print Your pass & bring it with you")

$p = $ppt.ActivePresentation

$oldText = "Print Your QR Code & bring it with you"
$newText = "Print Your pass & bring it with you"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $paraCount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $paraCount; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    if ($para.Text -eq $oldText) {
                        $runCount = $para.Runs().Count
                        if ($runCount -eq 1) {
                            # Single run holds the whole sentence - update it in place
                            # so formatting (font/size/spacing) is preserved exactly.
                            $run = $para.Runs(1, 1)
                            $run.Text = $newText
                        } else {
                            $para.Text = $newText
                        }
                    }
                }
            }
        }
    }
}
